$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AG (33): header date "11/11/21" mirrors the header style of AF1,
# and a value per person (rows 2-11) mirrors the plain style of the AF column.
# Values are forced to text (NumberFormat "@") before entry, matching the
# existing sheet's convention of storing every cell as a string, then the
# AF column's formatting is copied over (format-only paste) so the stored
# style index matches exactly (header keeps s="1", data cells keep default).

$headerCol = 32  # AF
$newCol = 33     # AG

$ws.Cells.Item(1, $newCol).NumberFormat = "@"
$ws.Cells.Item(1, $newCol).Value = "11/11/21"
$ws.Cells.Item(1, $headerCol).Copy()
$ws.Cells.Item(1, $newCol).PasteSpecial(-4122)

$values = @("0", "-1", "-2", "1", "0", "1", "0", "2", "0", "-1")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, $newCol).NumberFormat = "@"
    $ws.Cells.Item($row, $newCol).Value = $values[$i]
    $ws.Cells.Item($row, $headerCol).Copy()
    $ws.Cells.Item($row, $newCol).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
